$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = '28.951.16'
$ws.Range("E2").Value = '  +1.27%  '

# Row 3 - Ethereum
$ws.Range("D3").Value = '1.887.38'
$ws.Range("E3").Value = '  +0.74%  '

# Row 4 - TetherUSD
$ws.Range("E4").Value = '  -0.26%  '

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.19%  '

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4568'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.66%  '

# Row 8 - Cardano
$ws.Range("E8").Value = '  +1.46%  '

# Row 9 - Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07859'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.37%  '

# Row 10 - Polygon
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9861'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.54%  '

# Row 11 - Solana
$ws.Range("E11").Value = '  +1.86%  '

# Row 12 - WrappedEther
$ws.Range("D12").Value = '1.924.73'
$ws.Range("E12").Value = '  +1.41%  '

# Row 13 - Chainlink
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.032'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.48%  '

# Row 14 - Polkadot
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.693'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.93%  '

# Row 15 - TRON
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06951'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.31%  '

# Row 16 - Litecoin
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '88.09'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.63%  '

# Row 17 - BinanceUSD
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.004'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.20%  '

# Row 18 - ShibaInu
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009987'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.35%  '

# Row 19 - Avalanche
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.03'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.78%  '

# Row 20 - Dai
$ws.Range("E20").Value = '  -0.32%  '

# Row 21 - WrappedBTC
$ws.Range("D21").Value = '28.959.53'
$ws.Range("E21").Value = '  +1.27%  '

# Row 22 - Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.291'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.53%  '

# Row 23 - Cosmos
$ws.Range("E23").Value = '  +0.49%  '

# Row 24 - WrappedliquidstakedEther2.0
$ws.Range("D24").Value = '2.111.67'
$ws.Range("E24").Value = '  -0.28%  '

# Row 25 - Toncoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.053'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.25%  '

# Row 26 - Monero
$ws.Range("E26").Value = '  +0.90%  '

# Row 27 - EthereumClassic
$ws.Range("E27").Value = '  +0.98%  '

# Row 28 - InternetComputer(DFINITY)
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.963'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.03%  '

# Row 29 - LidoDAOToken
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.927'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.96%  '

# Row 30 - BitcoinCash
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '117.69'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.26%  '

# Row 31 - Stellar
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09333'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.94%  '

# Row 32 - ImmutableX
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9063'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.07%  '

# Row 33 - Filecoin
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.288'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.22%  '

# Row 34 - ARBITRUM
$ws.Range("E34").Value = '  +0.58%  '

# Row 35 - HuobiToken
$ws.Range("E35").Value = '  +0.25%  '

# Row 36 - TrustWalletToken
$ws.Range("E36").Value = '  +4.38%  '

# Row 37 - Hedera
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05767'
$ws.Range("D37").Style = "Normal"

# Row 38 - VeChain
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02072'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.47%  '

# Row 39 - Frax
$ws.Range("E39").Value = '  -0.29%  '

# Row 40 - TheSandbox(->FraxShare)
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.651'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.16%  '

# Row 41 - FraxShare(->TheSandbox)
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5670'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.70%  '

# Row 42 - Algorand
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1769'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.02%  '

# Row 43 - Aptos
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.734'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.90%  '

# Row 44 - RenderToken
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.262'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.98%  '

# Row 45 - EnergySwap
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '11.92'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.72%  '

# Row 46 - Decentraland
$ws.Range("E46").Value = '  +1.97%  '

# Row 47 - Cronos
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.07038'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.59%  '

# Row 48 - NEARProtocol
$ws.Range("E48").Value = '  +2.65%  '

# Row 49 - Quant
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '112.59'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.73%  '

# Row 50 - MXToken
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.515'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.16%  '

# Row 51 - WEMIXToken
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.078'
$ws.Range("D51").Style = "Normal"
